# Apply the work-log edit described by the commit:
#   "Fixed build systems, started ui, fixed average func, and implemented k-d tree"
#
# Concretely, on Sheet1 (the time-tracking log):
#   - Row 32 gets corrected: it now represents the next day (2019-03-14),
#     a new "End" time is recorded in column C, and the start time in
#     column B is corrected. The Hours formula (column D) recalculates
#     automatically from the shared formula C32-B32.
#   - The Total Hours sum in E3 recalculates automatically as a result.
#   - The view scrolls down a bit and the active selection moves to F32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32 data corrections -------------------------------------------------
# Day column (A32): was 3/13/2019 (43537), now 3/14/2019 (43538)
$ws.Range("A32").Value2 = 43538

# Start time (B32): corrected start time
$ws.Range("B32").Value2 = 0.41597222222222219

# End time (C32): newly filled in for this row
$ws.Range("C32").Value2 = 0.99930555555555556
# Match the time formatting used by the rest of this column (e.g. C31)
$ws.Range("C32").NumberFormat = $ws.Range("C31").NumberFormat

# Column D (Hours) keeps its existing shared formula (C32-B32); it will
# recalculate automatically, as will the Total Hours SUM in E3.

# --- View state ---------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F32").Select()
